$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.255.97'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '1.794.21'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '338.46'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4528'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +20.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3580'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +6.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '45.41'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.48%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.137'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07470'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.33%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.31'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.201'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.224'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.45%  '
$ws.Range('D16').Value = '1.791.81'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001081'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06687'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '81.03'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.18'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').Value = '28.214.37'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.83'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.86%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.381'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '20.35'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '153.39'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.377'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = '1.996.33'
$ws.Range('E29').Value = '  +1.66%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.267'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.86%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '132.13'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.073'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.09402'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +7.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02365'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '12.06'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.93%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6622'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06256'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.154'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2151'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.480'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.86%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.211'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.055'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9998'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.867'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.80%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.6052'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '128.35'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.57%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.019'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07079'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.158'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.93%  '
